$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear rows 22-27 entirely (table now stops at row 21)
$ws.Range("A22:K27").ClearContents()

# F4 and F11 previously held "Subtotals" labels; those rows are now
# regular Cost of Sales line items, so the hierarchy label is cleared.
$ws.Range("F4").ClearContents()
$ws.Range("F11").ClearContents()

# Row 2: Salaries- Operation Staff
$ws.Range("A2").Value = " Salaries- Operation Staff"
$ws.Range("B2").Value = 41704700
$ws.Range("C2").Value = 58420000
$ws.Range("D2").Value = $false
$ws.Range("E2").Value = "Cost of Sales"
$ws.Range("G2").Value = -16715300
$ws.Range("H2").Value = 12.6
$ws.Range("I2").Value = 12
$ws.Range("J2").Value = -28.6
$ws.Range("K2").Value = 0.6

# Row 3: Chicken
$ws.Range("A3").Value = " Chicken"
$ws.Range("B3").Value = 33505465
$ws.Range("C3").Value = 44669240
$ws.Range("D3").Value = $false
$ws.Range("E3").Value = "Cost of Sales"
$ws.Range("G3").Value = -11163775
$ws.Range("H3").Value = 10.1
$ws.Range("I3").Value = 9.2
$ws.Range("J3").Value = -25
$ws.Range("K3").Value = 1

# Row 4: Mutton
$ws.Range("A4").Value = " Mutton"
$ws.Range("B4").Value = 23844466
$ws.Range("C4").Value = 31502975
$ws.Range("D4").Value = $false
$ws.Range("E4").Value = "Cost of Sales"
$ws.Range("G4").Value = -7658509
$ws.Range("H4").Value = 7.2
$ws.Range("I4").Value = 6.5
$ws.Range("J4").Value = -24.3
$ws.Range("K4").Value = 0.7

# Row 5: Utilitites
$ws.Range("A5").Value = " Utilitites"
$ws.Range("B5").Value = 17459405.79
$ws.Range("C5").Value = 16942755
$ws.Range("D5").Value = $true
$ws.Range("E5").Value = "Cost of Sales"
$ws.Range("G5").Value = 516650.7899999991
$ws.Range("H5").Value = 5.3
$ws.Range("I5").Value = 3.5
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 1.8

# Row 6: Dairy & Bakery Items
$ws.Range("A6").Value = " Dairy & Bakery Items"
$ws.Range("B6").Value = 15684897
$ws.Range("C6").Value = 15563716
$ws.Range("D6").Value = $true
$ws.Range("E6").Value = "Cost of Sales"
$ws.Range("G6").Value = 121181
$ws.Range("H6").Value = 4.7
$ws.Range("I6").Value = 3.2
$ws.Range("J6").Value = 0.8
$ws.Range("K6").Value = 1.5

# Row 7: Other Dry Stock
$ws.Range("A7").Value = " Other Dry Stock"
$ws.Range("B7").Value = 15404163
$ws.Range("C7").Value = 14766946
$ws.Range("D7").Value = $true
$ws.Range("E7").Value = "Cost of Sales"
$ws.Range("G7").Value = 637217
$ws.Range("H7").Value = 4.6
$ws.Range("I7").Value = 3
$ws.Range("J7").Value = 4.3
$ws.Range("K7").Value = 1.6

# Row 8: Rents, rates & taxes
$ws.Range("A8").Value = " Rents, rates & taxes"
$ws.Range("B8").Value = 13199727
$ws.Range("C8").Value = 15480000
$ws.Range("D8").Value = $false
$ws.Range("E8").Value = "Cost of Sales"
$ws.Range("G8").Value = -2280273
$ws.Range("H8").Value = 4
$ws.Range("I8").Value = 3.2
$ws.Range("J8").Value = -14.7
$ws.Range("K8").Value = 0.8

# Row 9: Fuel & Transporation
$ws.Range("A9").Value = " Fuel & Transporation"
$ws.Range("B9").Value = 12481133
$ws.Range("C9").Value = 15052952
$ws.Range("D9").Value = $false
$ws.Range("E9").Value = "Cost of Sales"
$ws.Range("G9").Value = -2571819
$ws.Range("H9").Value = 3.8
$ws.Range("I9").Value = 3.1
$ws.Range("J9").Value = -17.1
$ws.Range("K9").Value = 0.7

# Row 10: Vegetables & Fruits
$ws.Range("A10").Value = " Vegetables & Fruits"
$ws.Range("B10").Value = 8824408
$ws.Range("C10").Value = 11670654
$ws.Range("D10").Value = $false
$ws.Range("E10").Value = "Cost of Sales"
$ws.Range("G10").Value = -2846246
$ws.Range("H10").Value = 2.7
$ws.Range("I10").Value = 2.4
$ws.Range("J10").Value = -24.4
$ws.Range("K10").Value = 0.3

# Row 11: Oil
$ws.Range("A11").Value = " Oil"
$ws.Range("B11").Value = 7329488
$ws.Range("C11").Value = 11911111
$ws.Range("D11").Value = $false
$ws.Range("E11").Value = "Cost of Sales"
$ws.Range("G11").Value = -4581623
$ws.Range("H11").Value = 2.2
$ws.Range("I11").Value = 2.4
$ws.Range("J11").Value = -38.5
$ws.Range("K11").Value = -0.2

# Row 12: Wages & Allowance
$ws.Range("A12").Value = " Wages & Allowance"
$ws.Range("B12").Value = 7320065
$ws.Range("C12").Value = 13365667
$ws.Range("D12").Value = $true
$ws.Range("E12").Value = "Cost of Sales"
$ws.Range("G12").Value = -6045602
$ws.Range("H12").Value = 2.2
$ws.Range("I12").Value = 2.7
$ws.Range("J12").Value = -45.2
$ws.Range("K12").Value = -0.5

# Row 13: Outsourced Food
$ws.Range("A13").Value = " Outsourced Food"
$ws.Range("B13").Value = 7064547
$ws.Range("C13").Value = 11473595
$ws.Range("D13").Value = $false
$ws.Range("E13").Value = "Cost of Sales"
$ws.Range("G13").Value = -4409048
$ws.Range("H13").Value = 2.1
$ws.Range("I13").Value = 2.4
$ws.Range("J13").Value = -38.4
$ws.Range("K13").Value = -0.2

# Row 14: Other Consumables
$ws.Range("A14").Value = " Other Consumables"
$ws.Range("B14").Value = 7011284.6
$ws.Range("C14").Value = 5096683
$ws.Range("D14").Value = $true
$ws.Range("E14").Value = "Cost of Sales"
$ws.Range("G14").Value = 1914601.6
$ws.Range("H14").Value = 2.1
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 37.6
$ws.Range("K14").Value = 1.1

# Row 15: Labour Out Sourced
$ws.Range("A15").Value = " Labour Out Sourced"
$ws.Range("B15").Value = 6882100
$ws.Range("C15").Value = 5827856
$ws.Range("D15").Value = $false
$ws.Range("E15").Value = "Cost of Sales"
$ws.Range("G15").Value = 1054244
$ws.Range("H15").Value = 2.1
$ws.Range("I15").Value = 1.2
$ws.Range("J15").Value = 18.1
$ws.Range("K15").Value = 0.9

# Row 16: Rice
$ws.Range("A16").Value = " Rice"
$ws.Range("B16").Value = 5963095
$ws.Range("C16").Value = 5116881
$ws.Range("D16").Value = $false
$ws.Range("E16").Value = "Cost of Sales"
$ws.Range("G16").Value = 846214
$ws.Range("H16").Value = 1.8
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 16.5
$ws.Range("K16").Value = 0.8

# Row 17: Beef
$ws.Range("A17").Value = " Beef"
$ws.Range("B17").Value = 5842129
$ws.Range("C17").Value = 5050412
$ws.Range("D17").Value = $false
$ws.Range("E17").Value = "Cost of Sales"
$ws.Range("G17").Value = 791717
$ws.Range("H17").Value = 1.8
$ws.Range("I17").Value = 1
$ws.Range("J17").Value = 15.7
$ws.Range("K17").Value = 0.7

# Row 18: Disposible Material
$ws.Range("A18").Value = " Disposible Material"
$ws.Range("B18").Value = 5742279
$ws.Range("C18").Value = 7628349
$ws.Range("D18").Value = $false
$ws.Range("E18").Value = "Cost of Sales"
$ws.Range("G18").Value = -1886070
$ws.Range("H18").Value = 1.7
$ws.Range("I18").Value = 1.6
$ws.Range("J18").Value = -24.7
$ws.Range("K18").Value = 0.2

# Row 19: Beverages & Soft drinks
$ws.Range("A19").Value = " Beverages & Soft drinks"
$ws.Range("B19").Value = 5032509
$ws.Range("C19").Value = 13011841
$ws.Range("D19").Value = $true
$ws.Range("E19").Value = "Cost of Sales"
$ws.Range("G19").Value = -7979332
$ws.Range("H19").Value = 1.5
$ws.Range("I19").Value = 2.7
$ws.Range("J19").Value = -61.3
$ws.Range("K19").Value = -1.1

# Row 20: Fish & Prawns
$ws.Range("A20").Value = " Fish & Prawns"
$ws.Range("B20").Value = 4122626
$ws.Range("C20").Value = 5899396
$ws.Range("D20").Value = $false
$ws.Range("E20").Value = "Cost of Sales"
$ws.Range("G20").Value = -1776770
$ws.Range("H20").Value = 1.2
$ws.Range("I20").Value = 1.2
$ws.Range("J20").Value = -30.1
$ws.Range("K20").Value = 0

# Row 21: Flour
$ws.Range("A21").Value = " Flour"
$ws.Range("B21").Value = 3858062
$ws.Range("C21").Value = 4155812
$ws.Range("D21").Value = $false
$ws.Range("E21").Value = "Cost of Sales"
$ws.Range("G21").Value = -297750
$ws.Range("H21").Value = 1.2
$ws.Range("I21").Value = 0.9
$ws.Range("J21").Value = -7.2
$ws.Range("K21").Value = 0.3
